$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily series. It belongs right
# before the existing row 78, so insert a blank row there (this shifts every
# row from 78 onward down by one, which is exactly what the target diff
# shows: old row 78 -> new row 79, old row 79 -> new row 80, ..., old row
# 196 -> new row 197).
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record's data.
$ws.Range("A78").Value = 9
$ws.Range("B78").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C78").Value = "Metropolitana"
$ws.Range("D78").Value = 44540
$ws.Range("E78").Value = 13
$ws.Range("F78").Value = 300000001
$ws.Range("G78").Value = "Rabanito"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 8800
$ws.Range("K78").Value = 2500
$ws.Range("L78").Value = 3000
$ws.Range("M78").Value = 2750
$ws.Range("N78").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O78").Value = "Provincia de Chacabuco"
$ws.Range("P78").Value = 28
$ws.Range("Q78").Value = 100
$ws.Range("R78").Value = "Hortaliza"
